$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style for the date column (column D), taken from an existing row
$dateFormat = $ws.Range("D190").NumberFormat

$newRows = @(
    @{ Row = 191; A = 6; B = "Mercado Mayorista Lo Valledor de Santiago"; C = "Metropolitana"; D = 45239; E = 13; F = 300000000; G = "Espárragos"; H = "Sin especificar"; I = "Banquete"; J = 480; K = 1400; L = 1500; M = 1448; N = "$/kilo"; O = "Provincia de Linares"; P = 1448; Q = 1; R = "Hortaliza" },
    @{ Row = 192; A = 6; B = "Mercado Mayorista Lo Valledor de Santiago"; C = "Metropolitana"; D = 45239; E = 13; F = 300000000; G = "Espárragos"; H = "Sin especificar"; I = "Primera"; J = 440; K = 1100; L = 1200; M = 1145; N = "$/kilo"; O = "Provincia de Linares"; P = 1145; Q = 1; R = "Hortaliza" },
    @{ Row = 193; A = 6; B = "Mercado Mayorista Lo Valledor de Santiago"; C = "Metropolitana"; D = 45239; E = 13; F = 300000000; G = "Espárragos"; H = "Sin especificar"; I = "Segunda"; J = 388; K = 900; L = 1000; M = 944; N = "$/kilo"; O = "Provincia de Linares"; P = 944; Q = 1; R = "Hortaliza" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value2 = $r.A
    $ws.Cells.Item($rowNum, 2).Value2 = $r.B
    $ws.Cells.Item($rowNum, 3).Value2 = $r.C
    $ws.Cells.Item($rowNum, 4).Value2 = $r.D
    $ws.Cells.Item($rowNum, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($rowNum, 5).Value2 = $r.E
    $ws.Cells.Item($rowNum, 6).Value2 = $r.F
    $ws.Cells.Item($rowNum, 7).Value2 = $r.G
    $ws.Cells.Item($rowNum, 8).Value2 = $r.H
    $ws.Cells.Item($rowNum, 9).Value2 = $r.I
    $ws.Cells.Item($rowNum, 10).Value2 = $r.J
    $ws.Cells.Item($rowNum, 11).Value2 = $r.K
    $ws.Cells.Item($rowNum, 12).Value2 = $r.L
    $ws.Cells.Item($rowNum, 13).Value2 = $r.M
    $ws.Cells.Item($rowNum, 14).Value2 = $r.N
    $ws.Cells.Item($rowNum, 15).Value2 = $r.O
    $ws.Cells.Item($rowNum, 16).Value2 = $r.P
    $ws.Cells.Item($rowNum, 17).Value2 = $r.Q
    $ws.Cells.Item($rowNum, 18).Value2 = $r.R
}

Write-Host "Added rows 191-193. UsedRange rows:" $ws.UsedRange.Rows.Count
